$wb = $excel.ActiveWorkbook
$ws1 = $wb.Worksheets.Item(1)
$ws3 = $wb.Worksheets.Item(3)

# Rename Data3 -> Login
$ws3.Name = "Login"

# Add Emp sheet right after Login
$empSheet = $wb.Worksheets.Add([System.Reflection.Missing]::Value, $ws3)
$empSheet.Name = "Emp"

# Copy the bordered style from Data1!A1:B1 onto Emp!A1:D1, then set values
$ws1.Range("A1:B1").Copy()
$empSheet.Range("A1:D1").PasteSpecial(-4122)

$empSheet.Range("A1").Value = "Aswini101"
$empSheet.Range("B1").Value = "Selenium"
$empSheet.Range("C1").Value = "Hyd"
$empSheet.Range("D1").Value = 101

$empSheet.Range("A1:D1").Select() | Out-Null

# Add User sheet right after Emp (left empty)
$userSheet = $wb.Worksheets.Add([System.Reflection.Missing]::Value, $empSheet)
$userSheet.Name = "User"

# Update the Login sheet's credentials and selection
$ws3.Range("A1").Value = "Admin1"
$ws3.Range("D5").Select() | Out-Null
